# Apply Sep-13-2023 cryptos data refresh to Sheet1 (columns D = Price, E = Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it looks like a number
# (e.g. "212.19"). Using a leading apostrophe forces Excel to store it as text,
# exactly like it would if a person typed it in, and resetting the style back to
# "Normal" afterwards avoids leaving a stray quote-prefixed / text-formatted cell style behind.
function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range("D2").Value = "26.224.45"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.605.19"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue "D5" "212.19"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  -0.16%  "
Set-TextValue "D9" "0.0615"
Set-TextValue "D10" "18.18"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.827.05"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.606.99"
$ws.Range("E13").Value = "  -0.32%  "
Set-TextValue "D14" "4.03"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "26.227.75"
$ws.Range("E16").Value = "  +0.32%  "
Set-TextValue "D17" "61.34"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  -0.13%  "
Set-TextValue "D20" "203.63"
$ws.Range("E20").Value = "  +2.34%  "
Set-TextValue "D21" "4.29"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +11.56%  "
Set-TextValue "D25" "144.67"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -6.78%  "
Set-TextValue "D28" "15.21"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.62%  "
Set-TextValue "D30" "0.0495"
$ws.Range("E30").Value = "  +4.07%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  +0.37%  "
Set-TextValue "D33" "2.93"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D36").Value = "1.138.63"
$ws.Range("E36").Value = "  +2.70%  "
Set-TextValue "D37" "0.0163"
$ws.Range("E37").Value = "  +6.23%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -0.58%  "
Set-TextValue "D40" "0.788"
$ws.Range("E40").Value = "  -0.44%  "
Set-TextValue "D41" "0.497"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D44").Value = "1.741.30"
$ws.Range("E44").Value = "  -0.45%  "
Set-TextValue "D45" "92.08"
$ws.Range("E45").Value = "  -1.19%  "
Set-TextValue "D46" "1.51"
$ws.Range("E46").Value = "  -3.46%  "
Set-TextValue "D47" "54.22"
$ws.Range("E47").Value = "  +0.55%  "
Set-TextValue "D48" "0.0507"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "0.0₇0948"
$ws.Range("E51").Value = "  -11.65%  "
